$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: NoTaxAddress billing/shipping dataset
$ws.Range("A20").Value = "NoTaxAddress"
$ws.Range("F20").Value = "Test"
$ws.Range("N20").Value = "879 Killens Pond Rd"
$ws.Range("O20").Value = "Harrington"
$ws.Range("P20").Value = "United States"
$ws.Range("Q20").Value = "Delaware"
$ws.Range("R20").Formula = "'19952"
$ws.Range("G20").Value = "qa"
$ws.Range("S20").Value = 9898989898

# Row 21: BillingDetails dataset (different billing and shipping)
$ws.Range("A21").Value = "BillingDetails"
$ws.Range("F21").Value = "QA"
$ws.Range("G21").Value = "TEST"
$ws.Range("N21").Value = "6 Walnut Valley Dr"
$ws.Range("O21").Value = "Little Rock"
$ws.Range("P21").Value = "United States"
$ws.Range("Q21").Value = "Arkansas"
$ws.Range("R21").Formula = "'72211"
$ws.Range("S21").Value = 9898989898

# Update selection to match the authored state
[void]$ws.Range("D20").Select()
